# Insert two new weekly price-report rows (Espinaca, Mercado Mayorista Lo
# Valledor de Santiago) at the top of the date-ordered data block, pushing
# the existing rows 247-356 down to 249-358.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("247:248").Insert()

# New row 247 - "Provincia de Chacabuco" origin, calidad Primera
$ws.Range("A247").Value2 = 6
$ws.Range("B247").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C247").Value2 = "Metropolitana"
$ws.Range("D247").Value2 = 44466
$ws.Range("E247").Value2 = 13
$ws.Range("F247").Value2 = 100112012
$ws.Range("G247").Value2 = "Espinaca"
$ws.Range("H247").Value2 = "Sin especificar"
$ws.Range("I247").Value2 = "Primera"
$ws.Range("J247").Value2 = 380
$ws.Range("K247").Value2 = 5500
$ws.Range("L247").Value2 = 6000
$ws.Range("M247").Value2 = 5724
$ws.Range("N247").Value2 = "`$/cuna 10 kilos"
$ws.Range("O247").Value2 = "Provincia de Chacabuco"
$ws.Range("P247").Value2 = 572
$ws.Range("Q247").Value2 = 10
$ws.Range("R247").Value2 = "Hortaliza"

# New row 248 - "Región Metropolitana" origin, calidad Primera
$ws.Range("A248").Value2 = 6
$ws.Range("B248").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C248").Value2 = "Metropolitana"
$ws.Range("D248").Value2 = 44466
$ws.Range("E248").Value2 = 13
$ws.Range("F248").Value2 = 100112012
$ws.Range("G248").Value2 = "Espinaca"
$ws.Range("H248").Value2 = "Sin especificar"
$ws.Range("I248").Value2 = "Primera"
$ws.Range("J248").Value2 = 430
$ws.Range("K248").Value2 = 5500
$ws.Range("L248").Value2 = 6000
$ws.Range("M248").Value2 = 5709
$ws.Range("N248").Value2 = "`$/cuna 10 kilos"
$ws.Range("O248").Value2 = "Región Metropolitana"
$ws.Range("P248").Value2 = 571
$ws.Range("Q248").Value2 = 10
$ws.Range("R248").Value2 = "Hortaliza"
